# "update on 15 nov" -- add new GAMMA / THETA / VEGA / strategies notes
# below the existing "options greeks" table (rows 38-62), plus a small
# "mtm" note in K34 and a new column-C width tweak.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new content, entered in the same order the author must have typed it
#     (this keeps the shared-string table's append order identical to the
#     authored workbook) ---

# GAMMA block
$ws.Range("A38").Value = "GAMMA"
$ws.Range("A38").Font.Bold = $true
$ws.Range("B38").Value = "gamma tells how much delta will change when stock price  changes"
$ws.Range("B39").Value = "delta is speed and gamma is acceltrator "
$ws.Range("B40").Value = "ATM has highest Gamma"
$ws.Range("B41").Value = "gamma increase to close of expiration"

# THETA header, plus a stray "mtm" note added over in K34 at the same time
$ws.Range("A43").Value = "THETA"
$ws.Range("A43").Font.Bold = $true
$ws.Range("K34").Value = "mtm"
$ws.Range("B43").Value = "time decays"
$ws.Range("B44").Value = "at expiry theta reaches to 0 "
$ws.Range("B45").Value = "ATM has the highest theta"

# VEGA block
$ws.Range("A48").Value = "VEGA"
$ws.Range("A48").Font.Bold = $true
$ws.Range("B48").Value = "volatitllity"
$ws.Range("B49").Value = "valatility make large part of extrinsic value"
$ws.Range("B50").Value = "buy option at low volatility"
$ws.Range("B51").Value = "sell option at high volatility"

# back-fill B46 under THETA
$ws.Range("B46").Value = "neg-ve for buy and +ve for sale"

# implied volatility notes
$ws.Range("A53").Value = "implied volitility"
$ws.Range("B54").Value = "india vix is IV"

# strategies list
$ws.Range("A58").Value = "Stretegies"
$ws.Range("A59").Value = "1. long call"
$ws.Range("A60").Value = "2. sell  put"
$ws.Range("B59").Value = "don" + [char]0x2019 + "t buy on high IV"
$ws.Range("B60").Value = "don" + [char]0x2019 + "t sell on market crash"
$ws.Range("A62").Value = "3. Bull call spread"

# --- small formatting touch-ups ---

# new column C width (stored width 10 == ColumnWidth 10 - 5/6)
$ws.Columns("C:C").ColumnWidth = 9.166666666666666

# hide gridlines and leave the selection where the author left it
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A63").Select() | Out-Null
